$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MergeSort")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
